$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quantum-state column headers: now reporting the swap-test ancilla outcomes ---
# |0> -> |1>, |1> -> |01>  (D1 ends up showing |01>, E1 ends up showing |1>)
$ws.Range("D1").Value = "|01>"
$ws.Range("E1").Value = "|1>"

# --- Updated measured probability (D2) and the dependent formulas in column F ---
# Swap-test estimator: F = 3 - 4*D  (was F = 2*D - 1)
$ws.Range("D2").Value = 0.466

$ws.Range("F2").Formula = "=3-4*D2"
$ws.Range("F3").Formula = "=3-4*D3"
$ws.Range("F4").Formula = "=3-4*D4"
$ws.Range("F5").Formula = "=3-4*D5"
$ws.Range("F6").Formula = "=3-4*D6"
$ws.Range("F7").Formula = "=3-4*D7"
$ws.Range("F8").Formula = "=3-4*D8"
$ws.Range("F9").Formula = "=3-4*D9"

# --- Remove the old single-qubit circuit-diagram picture (swapped out for the swap test) ---
if ($ws.Shapes.Count -gt 0) {
    $ws.Shapes.Item(1).Delete()
}

# --- Restore the window position / move the saved selection cursor to D3 ---
$win = $excel.Windows.Item(1)
$win.Left = 6855
$win.Top = 2175

$ws.Range("D3").Select()
